$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.04027
$ws.Range("H2").Value = 0.12081
$ws.Range("I2").Value = 0.01318991723029425
$ws.Range("J2").Value = 0.01318991723029425
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 32.63563666666667
$ws.Range("N2").Value = 97.90691000000001
$ws.Range("O2").Value = 0.9900101876891448
$ws.Range("P2").Value = 0.9900101876891446
$ws.Range("Q2").Value = 1.314237088566667
$ws.Range("R2").Value = 11.8281337971
$ws.Range("S2").Value = 0.0130581524327679
$ws.Range("T2").Value = 0.01305815243276789

$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.04027
$ws.Range("H3").Value = 0.12081
$ws.Range("I3").Value = 0.01318991723029425
$ws.Range("J3").Value = 0.01318991723029425
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.05920466666666666
$ws.Range("N3").Value = 0.177614
$ws.Range("O3").Value = 0.001795988347259859
$ws.Range("P3").Value = 0.001795988347259859
$ws.Range("Q3").Value = 0.002384171926666667
$ws.Range("R3").Value = 0.02145754734
$ws.Range("S3").Value = 0.00002368893764693051
$ws.Range("T3").Value = 0.0000236889376469305

$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.04027
$ws.Range("H4").Value = 0.12081
$ws.Range("I4").Value = 0.01318991723029425
$ws.Range("J4").Value = 0.01318991723029425
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.270109
$ws.Range("N4").Value = 0.810327
$ws.Range("O4").Value = 0.008193823963595435
$ws.Range("P4").Value = 0.008193823963595434
$ws.Range("Q4").Value = 0.01087728943
$ws.Range("R4").Value = 0.09789560487
$ws.Range("S4").Value = 0.0001080758598794254
$ws.Range("T4").Value = 0.0001080758598794254

$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.9943730000000001
$ws.Range("H5").Value = 2.983119
$ws.Range("I5").Value = 0.3256940046198011
$ws.Range("J5").Value = 0.325694004619801
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 32.63563666666667
$ws.Range("N5").Value = 97.90691000000001
$ws.Range("O5").Value = 0.9900101876891448
$ws.Range("P5").Value = 0.9900101876891446
$ws.Range("Q5").Value = 32.45199593914334
$ws.Range("R5").Value = 292.0679634522901
$ws.Range("S5").Value = 0.3224403826428785
$ws.Range("T5").Value = 0.3224403826428783

$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.9943730000000001
$ws.Range("H6").Value = 2.983119
$ws.Range("I6").Value = 0.3256940046198011
$ws.Range("J6").Value = 0.325694004619801
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.05920466666666666
$ws.Range("N6").Value = 0.177614
$ws.Range("O6").Value = 0.001795988347259859
$ws.Range("P6").Value = 0.001795988347259859
$ws.Range("Q6").Value = 0.05887152200733333
$ws.Range("R6").Value = 0.5298436980660001
$ws.Range("S6").Value = 0.0005849426370695614
$ws.Range("T6").Value = 0.0005849426370695611

$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.9943730000000001
$ws.Range("H7").Value = 2.983119
$ws.Range("I7").Value = 0.3256940046198011
$ws.Range("J7").Value = 0.325694004619801
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.270109
$ws.Range("N7").Value = 0.810327
$ws.Range("O7").Value = 0.008193823963595435
$ws.Range("P7").Value = 0.008193823963595434
$ws.Range("Q7").Value = 0.268589096657
$ws.Range("R7").Value = 2.417301869913
$ws.Range("S7").Value = 0.002668679339853089
$ws.Range("T7").Value = 0.002668679339853088

$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 2.018446666666666
$ws.Range("H8").Value = 6.055339999999999
$ws.Range("I8").Value = 0.6611160781499047
$ws.Range("J8").Value = 0.6611160781499047
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 32.63563666666667
$ws.Range("N8").Value = 97.90691000000001
$ws.Range("O8").Value = 0.9900101876891448
$ws.Range("P8").Value = 0.9900101876891446
$ws.Range("Q8").Value = 65.87329204437778
$ws.Range("R8").Value = 592.8596283994
$ws.Range("S8").Value = 0.6545116526134985
$ws.Range("T8").Value = 0.6545116526134983

$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 2.018446666666666
$ws.Range("H9").Value = 6.055339999999999
$ws.Range("I9").Value = 0.6611160781499047
$ws.Range("J9").Value = 0.6611160781499047
$ws.Range("K9").Value = 1
$ws.Range("L9").Value = 0.3333333333333333
$ws.Range("M9").Value = 0.05920466666666666
$ws.Range("N9").Value = 0.177614
$ws.Range("O9").Value = 0.001795988347259859
$ws.Range("P9").Value = 0.001795988347259859
$ws.Range("Q9").Value = 0.1195014620844444
$ws.Range("R9").Value = 1.07551315876
$ws.Range("S9").Value = 0.001187356772543367
$ws.Range("T9").Value = 0.001187356772543367

$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 2.018446666666666
$ws.Range("H10").Value = 6.055339999999999
$ws.Range("I10").Value = 0.6611160781499047
$ws.Range("J10").Value = 0.6611160781499047
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0.270109
$ws.Range("N10").Value = 0.810327
$ws.Range("O10").Value = 0.008193823963595435
$ws.Range("P10").Value = 0.008193823963595434
$ws.Range("Q10").Value = 0.5452006106866666
$ws.Range("R10").Value = 4.90680549618
$ws.Range("S10").Value = 0.005417068763862921
$ws.Range("T10").Value = 0.005417068763862921
